$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for Albahaca (Femacal de La Calera).
# It belongs right after the existing row 77, so insert a fresh row at 78
# and push every following record down by one (row 180 -> 181).
$ws.Rows("78").Insert()

# Populate the newly inserted row 78 with the new observation's data.
$ws.Range("A78").Value = 3
$ws.Range("B78").Value = "Femacal de La Calera"
$ws.Range("C78").Value = "Coquimbo"
$ws.Range("D78").Value = 44893
$ws.Range("E78").Value = 5
$ws.Range("F78").Value = 100112052
$ws.Range("G78").Value = "Albahaca"
$ws.Range("H78").Value = "Sin especificar"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 135
$ws.Range("K78").Value = 6000
$ws.Range("L78").Value = 6500
$ws.Range("M78").Value = 6241
$ws.Range("N78").Value = "$/docena de matas"
$ws.Range("O78").Value = "Provincia de Quillota"
$ws.Range("P78").Value = 1040
$ws.Range("Q78").Value = 6
$ws.Range("R78").Value = "Hortaliza"
